$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# A new task ("Encontrar bebedouros perto de Moscavide") was inserted as row 11
# of the Gantt table. Rows 11-20 (the rest of "PRIMEIRA ENTREGA" through
# "TERCEIRA ENTREGA") need to shift down one row to 12-21. Rows below the
# table (B26/B27 - "inicio"/"fim") are NOT part of this contiguous block and
# stay where they are, so we shift the B:E block by copying row-by-row from
# the bottom up instead of doing a real "insert entire row" (which would
# also push the unrelated summary rows down).
for ($r = 20; $r -ge 11; $r--) {
  $srcRow = $ws.Range("B" + $r + ":E" + $r)
  $dstRow = $ws.Range("B" + ($r + 1) + ":E" + ($r + 1))
  $srcRow.Copy($dstRow)
  # Re-key the duration formula to the row it now lives in.
  $ws.Range("E" + ($r + 1)).Formula = "=C" + ($r + 1) + "+D" + ($r + 1)
}
$excel.CutCopyMode = 0

# Fill the freed-up row 11 with the new task.
$ws.Range("B11").Value = "Encontrar bebedouros perto de Moscavide"
$ws.Range("C11").Value = 45928
$ws.Range("D11").Value = 5
$ws.Range("E11").Formula = "=C11+D11"

# "fim" (C27) tracked the end of the TERCEIRA ENTREGA row, which used to be
# row 19 and is now row 20.
$ws.Range("C27").Formula = "=E20"

# Refresh the sheet view: scrolled down to the summary rows, zoomed to 77%,
# with I14 as the active selection.
$ws.Activate()
$excel.ActiveWindow.Zoom = 77
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I14").Select()
